$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture existing values for rows 2..11, columns B..G before making any changes,
# since each row's new B:F values come from the row above's old values (a shift),
# and new row 2 data is freshly supplied.
$oldValues = @{}
for ($r = 2; $r -le 11; $r++) {
    $oldValues[$r] = @(
        $ws.Cells.Item($r, 2).Value(),  # B
        $ws.Cells.Item($r, 3).Value(),  # C
        $ws.Cells.Item($r, 4).Value(),  # D
        $ws.Cells.Item($r, 5).Value(),  # E
        $ws.Cells.Item($r, 6).Value(),  # F
        $ws.Cells.Item($r, 7).Value()   # G
    )
}

# New freshly-computed values for row 2 (B:F)
$newRow2 = @(0.1184146901506048, 0.5185973060008381, 0.4691159408045538, 0.6849203901217672, 0.698284046247675)

# Write row 2: new B:F values, G = old row2 G + 1
$ws.Cells.Item(2, 2).Value = $newRow2[0]
$ws.Cells.Item(2, 3).Value = $newRow2[1]
$ws.Cells.Item(2, 4).Value = $newRow2[2]
$ws.Cells.Item(2, 5).Value = $newRow2[3]
$ws.Cells.Item(2, 6).Value = $newRow2[4]
$ws.Cells.Item(2, 7).Value = $oldValues[2][5] + 1

# Write rows 3..11: B:F shift down from the row above's old values, G = old G + 1
for ($r = 3; $r -le 11; $r++) {
    $prev = $oldValues[$r - 1]
    $cur = $oldValues[$r]
    $ws.Cells.Item($r, 2).Value = $prev[0]
    $ws.Cells.Item($r, 3).Value = $prev[1]
    $ws.Cells.Item($r, 4).Value = $prev[2]
    $ws.Cells.Item($r, 5).Value = $prev[3]
    $ws.Cells.Item($r, 6).Value = $prev[4]
    $ws.Cells.Item($r, 7).Value = $cur[5] + 1
}
